$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49 (pushes current rows 49-57 down to 50-58),
# copying formatting (incl. the date-format style on column D) from row 49.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly price record.
$ws.Cells.Item(49, 1).Value = 1
$ws.Cells.Item(49, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value = 44543
$ws.Cells.Item(49, 5).Value = 15
$ws.Cells.Item(49, 6).Value = 100112038
$ws.Cells.Item(49, 7).Value = "Cebollín baby"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 250
$ws.Cells.Item(49, 11).Value = 2500
$ws.Cells.Item(49, 12).Value = 3000
$ws.Cells.Item(49, 13).Value = 2750
$ws.Cells.Item(49, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(49, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value = 1375
$ws.Cells.Item(49, 17).Value = 2
$ws.Cells.Item(49, 18).Value = "Hortaliza"
